$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 17:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 851586
$ws.Range("C4").Value = 2869
$ws.Range("D4").Value = 84117
$ws.Range("E4").Value = 719661
$ws.Range("F4").Value = 14334
$ws.Range("G4").Value = 149
$ws.Range("H4").Value = 47808

# Row 16 - Canada
$ws.Range("B16").Value = 40824
$ws.Range("C16").Value = 634
$ws.Range("E16").Value = 24810
$ws.Range("G16").Value = 54
$ws.Range("H16").Value = 2028

# Row 28 - Chile
$ws.Range("B28").Value = 11812
$ws.Range("C28").Value = 516
$ws.Range("D28").Value = 5804
$ws.Range("E28").Value = 5840
$ws.Range("F28").Value = 411
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = 168

# Row 48 - Republica Dominicana
$ws.Range("B48").Value = 5543
$ws.Range("C48").Value = 243
$ws.Range("E48").Value = 4697
$ws.Range("G48").Value = 5
$ws.Range("H48").Value = 265

# Row 61 - Grecia
$ws.Range("B61").Value = 2463
$ws.Range("C61").Value = 55
$ws.Range("E61").Value = 1761
$ws.Range("F61").Value = 52

# Row 68 - Uzbekistan
$ws.Range("D68").Value = 561
$ws.Range("E68").Value = 1167

# Row 70 - Irak
$ws.Range("B70").Value = 1677
$ws.Range("C70").Value = 46
$ws.Range("D70").Value = 1171
$ws.Range("E70").Value = 423

# Row 90 - Republica de Chipre
$ws.Range("B90").Value = 795
$ws.Range("C90").Value = 5
$ws.Range("E90").Value = 684

# Row 118 - Montenegro
$ws.Range("D118").Value = 123
$ws.Range("E118").Value = 188
